# "Edit Heatmap after F/G reorganization"
# - Correct the trans-distance value in row 65 of the "Freq. of Occurrence"
#   sheet (I65: 89.24 -> 98.24), which shifted after the F/G column
#   reorganization.
# - Leave the sheet scrolled/selected where the author was working
#   (row ~50 in view, J69 selected) when they saved the file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Freq. of Occurrence")

# Data correction
$ws.Range("I65").Value = 98.24

# Restore the view/selection state as it was left after the edit
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 50
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J69").Select()
